$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" (strike count -> K) values for column G, rows 2-16.
$kValues = @{
    2  = 5
    3  = 2
    4  = 4
    5  = 5
    6  = 0
    7  = 0
    8  = 3
    9  = 3
    10 = 1
    11 = 4
    12 = 7
    13 = 2
    14 = 0
    15 = 1
    16 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
